$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Ray:" row (row 2) results for each benchmark column
# (order matches shared-string table insertion order in the target workbook)
$ws.Range("C2").Value = "35.44 seconds"
$ws.Range("B2").Value = "33.04 seconds"
$ws.Range("D2").Value = "32.88 seconds"

# Update the active selection to D2 as in the edited workbook
$ws.Range("D2").Select()
